# Refresh cryptocurrency price/volume snapshot (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.752.44"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").Value = "3.452.19"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'582.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").Value = "'147.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.11%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.477"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "'7.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "4.038.87"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("D13").Value = "'29.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.24%  "
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "3.451.06"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "62.818.87"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "'6.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").Value = "'14.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.50%  "
$ws.Range("D20").Value = "'9.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.06%  "
$ws.Range("D21").Value = "'396.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.24%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "'0.564"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'75.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'0.0000118"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.71%  "
$ws.Range("D26").Value = "3.574.31"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "'0.190"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'7.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.41%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'8.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.84%  "
$ws.Range("D31").Value = "'2.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "'1.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.12%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'23.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").Value = "'5.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.14%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.51%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'7.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").Value = "'169.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "3.485.40"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'29.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.48%  "
$ws.Range("D41").Value = "'0.0769"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("D43").Value = "'4.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("E44").Value = "  +4.97%  "
$ws.Range("E45").Value = "  +7.30%  "
$ws.Range("D46").Value = "2.533.38"
$ws.Range("E46").Value = "  +4.36%  "
$ws.Range("D47").Value = "'23.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("D48").Value = "'6.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "'2.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.18%  "
